$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The existing "_GoBack" bookmark (sitting next to "Cuadro de
#    precios acordados:") is stale - it marks a prior edit location.
#    Word recreates "_GoBack" at the *new* last-edit spot, so drop the
#    old one here; it gets re-added below once we know where the last
#    edit actually lands.
# ------------------------------------------------------------------
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
}

# ------------------------------------------------------------------
# 2) First price: "$ 3.300,00" -> "$ 3.500,00"
# ------------------------------------------------------------------
$d.Content.Find.Execute("`$ 3.300,00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "`$ 3.500,00", 2)

# ------------------------------------------------------------------
# 3) Second price: "$ 3.200,00" -> "$ 3.400,00", and this is where the
#    author's cursor was left, so "_GoBack" is re-created right after
#    "$ 3.4" (before the trailing "00,00").
# ------------------------------------------------------------------
$found = $d.Content
$found.Find.Execute("`$ 3.200,00", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$priceStart = $found.Start
$priceEnd = $found.End
$found = $null

$priceRange = $d.Range($priceStart, $priceEnd)
$priceRange.Text = "`$ 3.400,00"
$priceRange = $null

$goBackPos = $priceStart + 5
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
